$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Implemented delete for single users — fill in the newly-added error-code
# rows (10036, 10037, 10038, 10051) with their message key (column B) and
# severity level (column D). Shared strings are appended in the same order
# Excel would have written them (row 39, then 38, then 53, then 40) so the
# resulting sharedStrings.xml ordering lines up with the real commit.

$ws.Range("B39").Value = "message_10037_user_cannot_delete_own_record"
$ws.Range("D39").Value = "Error"

$ws.Range("B38").Value = "message_10036_password_updated_successfully"
$ws.Range("D38").Value = "Success"

$ws.Range("B53").Value = "message_10051_email_address_format_invalid"
$ws.Range("D53").Value = "Error"

$ws.Range("B40").Value = "message_10038_username_not_available"
$ws.Range("D40").Value = "Error"

# Mirror the author's final selection in the sheet view (the scroll-position
# "topLeftCell" is transient UI state that this host doesn't persist, but the
# active-cell selection is).
$ws.Range("B40").Select()
